$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.816.13"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.90%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.807.96"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.27%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.43%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "421.59"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +4.06%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.76"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.28%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.801.97"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +4.14%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.719"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.48%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.160"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.50%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000346"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +13.71%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.64"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.83%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.418.78"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.80%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "10.13"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.53"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +16.66%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.50%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.811.45"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.48%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.64"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "66.983.25"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.08"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.44%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "406.18"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.25"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.69%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.78"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.24%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.84%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "37.05"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.97"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +6.45%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.54%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.45"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +6.15%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.03"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +30.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "720.05"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +7.46%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.70"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.74%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.62%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.72%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.06%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -5.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "38.50"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.17"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.41"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +24.52%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0750"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +19.30%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.90%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.91"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.28%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.13%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.44%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.33"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.74%  "

$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.12"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.40%  "

$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "143.11"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.42%  "

$ws.Range("B48").Value = "TheGraph"
$ws.Range("C48").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.310"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +7.18%  "

$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.03"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.80"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "25.57"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -12.97%  "
